$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the lapso (date range) header text
$ws.Range("A1").Value = "Lapso: 2022-10-24 al 2022-11-22"

# Update rows 3-5 with the new (alphabetically sorted) specialist order and counts
$ws.Range("A3").Value = "Alberto Chinsky"
$ws.Range("B3").Value = 3

$ws.Range("A4").Value = "Selene Montaño"
$ws.Range("B4").Value = 3

$ws.Range("A5").Value = "Stefania Beatriz Marco"
$ws.Range("B5").Value = 3
